$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 618
$ws.Range("A618").Value = 0
$ws.Range("B618").Value = 20082200
$ws.Range("C618").Value = 0
$ws.Range("D618").Value = 0
$ws.Range("E618").Value = 0
$ws.Range("F618").Value = 0

# Row 619
$ws.Range("A619").Value = 0
$ws.Range("B619").Value = 20082300
$ws.Range("C619").Value = 0
$ws.Range("D619").Value = 0
$ws.Range("E619").Value = 0
$ws.Range("F619").Value = 0

# Row 620
$ws.Range("A620").Value = "9127963B3"
$ws.Range("B620").Value = 20082400
$ws.Range("C620").Value = 105000000000
$ws.Range("D620").Value = 54000000000
$ws.Range("E620").Value = 0.5142857142857142
$ws.Range("F620").Value = 20082700

# Row 621
$ws.Range("A621").Value = "912796XE4"
$ws.Range("B621").Value = 20082400
$ws.Range("C621").Value = 105000000000
$ws.Range("D621").Value = 51000000000
$ws.Range("E621").Value = 0.4857142857142857
$ws.Range("F621").Value = 20082700

# Row 622
$ws.Range("A622").Value = "91282CAG6"
$ws.Range("B622").Value = 20082500
$ws.Range("C622").Value = 110000000000
$ws.Range("D622").Value = 50000000000
$ws.Range("E622").Value = 0.4545454545454545
$ws.Range("F622").Value = 20083100

# Row 623
$ws.Range("A623").Value = "912796TN9"
$ws.Range("B623").Value = 20082500
$ws.Range("C623").Value = 110000000000
$ws.Range("D623").Value = 30000000000
$ws.Range("E623").Value = 0.2727272727272727
$ws.Range("F623").Value = 20082700

# Row 624
$ws.Range("A624").Value = "9127963L1"
$ws.Range("B624").Value = 20082500
$ws.Range("C624").Value = 110000000000
$ws.Range("D624").Value = 30000000000
$ws.Range("E624").Value = 0.2727272727272727
$ws.Range("F624").Value = 20082700

# Row 625
$ws.Range("A625").Value = "9127965C9"
$ws.Range("B625").Value2 = "'20082600"
$ws.Range("C625").Value = 128000000000
$ws.Range("D625").Value2 = "'25000000000"
$ws.Range("E625").Value = 0.1953125
$ws.Range("F625").Value2 = "'20090100"

# Row 626
$ws.Range("A626").Value = "912796B65"
$ws.Range("B626").Value2 = "'20082600"
$ws.Range("C626").Value = 128000000000
$ws.Range("D626").Value2 = "'30000000000"
$ws.Range("E626").Value = 0.234375
$ws.Range("F626").Value2 = "'20090100"

# Row 627
$ws.Range("A627").Value = "91282CAJ0"
$ws.Range("B627").Value2 = "'20082600"
$ws.Range("C627").Value = 128000000000
$ws.Range("D627").Value2 = "'51000000000"
$ws.Range("E627").Value = 0.3984375
$ws.Range("F627").Value2 = "'20083100"

# Row 628
$ws.Range("A628").Value = "91282CAA9"
$ws.Range("B628").Value2 = "'20082600"
$ws.Range("C628").Value = 128000000000
$ws.Range("D628").Value2 = "'22000000000"
$ws.Range("E628").Value = 0.171875
$ws.Range("F628").Value2 = "'20082800"

# Row 629
$ws.Range("A629").Value = "9127964A4"
$ws.Range("B629").Value2 = "'20082700"
$ws.Range("C629").Value = 112000000000
$ws.Range("D629").Value2 = "'30000000000"
$ws.Range("E629").Value = 0.2678571428571428
$ws.Range("F629").Value2 = "'20090100"

# Row 630
$ws.Range("A630").Value = "9127964K2"
$ws.Range("B630").Value2 = "'20082700"
$ws.Range("C630").Value = 112000000000
$ws.Range("D630").Value2 = "'35000000000"
$ws.Range("E630").Value = 0.3125
$ws.Range("F630").Value2 = "'20090100"

# Row 631
$ws.Range("A631").Value = "91282CAH4"
$ws.Range("B631").Value2 = "'20082700"
$ws.Range("C631").Value = 112000000000
$ws.Range("D631").Value2 = "'47000000000"
$ws.Range("E631").Value = 0.4196428571428572
$ws.Range("F631").Value2 = "'20083100"

# Row 632
$ws.Range("A632").Value = 0
$ws.Range("B632").Value2 = "'20082800"
$ws.Range("C632").Value = 0
$ws.Range("D632").Value = 0
$ws.Range("E632").Value = 0
$ws.Range("F632").Value = 0

# Row 633
$ws.Range("A633").Value = 0
$ws.Range("B633").Value2 = "'20082900"
$ws.Range("C633").Value = 0
$ws.Range("D633").Value = 0
$ws.Range("E633").Value = 0
$ws.Range("F633").Value = 0

# Row 634
$ws.Range("A634").Value = 0
$ws.Range("B634").Value2 = "'20083000"
$ws.Range("C634").Value = 0
$ws.Range("D634").Value = 0
$ws.Range("E634").Value = 0
$ws.Range("F634").Value = 0

# Row 635
$ws.Range("A635").Value = 0
$ws.Range("B635").Value2 = "'20083100"
$ws.Range("C635").Value = 0
$ws.Range("D635").Value = 0
$ws.Range("E635").Value = 0
$ws.Range("F635").Value = 0
